$wb = $excel.ActiveWorkbook

# --- Sheet "catalogo": fill in missing photo filenames for the Paris gorro rows ---
$catalogo = $wb.Worksheets.Item("catalogo")

$catalogo.Range("E35").Value = "paris negro.jpg"
$catalogo.Range("E36").Value = "paris marron.jpg"
$catalogo.Range("E37").Value = "paris crudo.jpg"

# Fill in missing categoría value for the "Granda" row
$catalogo.Range("C43").Value = "Granada"

# --- Sheet "datos": categoría helper table - move "Granada" entry up to the
# first empty slot (row 9) and clear the old trailing occurrence (row 15) ---
$datos = $wb.Worksheets.Item("datos")

$datos.Range("C9").Value = "Granada"
$datos.Range("C15").ClearContents()
